$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 253.92857
$ws.Range("I5").Value = 334.66666
$ws.Range("J5").Value = 108.6
$ws.Range("K5").Value = 334.66666
$ws.Range("L5").Value = 108.6
$ws.Range("M5").Value = -219.66666
$ws.Range("N5").Value = -338.6

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2386572
$ws.Range("J17").Value = 2386572
$ws.Range("L17").Value = 7159716
$ws.Range("N17").Value = -7160052

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 10140.454
$ws.Range("I55").Value = 1760.8334
$ws.Range("J55").Value = 20196
$ws.Range("K55").Value = 1760.8334
$ws.Range("L55").Value = 20196
$ws.Range("M55").Value = -1546.8334
$ws.Range("N55").Value = -20624

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H61").Value = 55
$ws.Range("I61").Value = 55
$ws.Range("K61").Value = 165
$ws.Range("M61").Value = 7

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 3374.077
$ws.Range("I64").Value = 2930.6897
$ws.Range("J64").Value = 4659.9
$ws.Range("K64").Value = 2930.6897
$ws.Range("L64").Value = 4659.9
$ws.Range("M64").Value = -2682.6897
$ws.Range("N64").Value = -5155.9

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 3374.077
$ws.Range("I67").Value = 2930.6897
$ws.Range("J67").Value = 4659.9
$ws.Range("K67").Value = 2930.6897
$ws.Range("L67").Value = 4659.9
$ws.Range("M67").Value = -2072.6897
$ws.Range("N67").Value = -6375.9

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H74").Value = 3240.0334
$ws.Range("I74").Value = 2793.923
$ws.Range("J74").Value = 3581.1765
$ws.Range("K74").Value = 2793.923
$ws.Range("L74").Value = 3581.1765
$ws.Range("M74").Value = -1857.923
$ws.Range("N74").Value = -5453.1765

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H77").Value = 3240.0334
$ws.Range("I77").Value = 2793.923
$ws.Range("J77").Value = 3581.1765
$ws.Range("K77").Value = 13969.615
$ws.Range("L77").Value = 17905.8825
$ws.Range("M77").Value = -9289.614999999998
$ws.Range("N77").Value = -27265.8825

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 324.57144
$ws.Range("I92").Value = 253.47368
$ws.Range("K92").Value = 253.47368
$ws.Range("M92").Value = 994.5263199999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H100").Value = 20898448
$ws.Range("I100").Value = 33335154
$ws.Range("K100").Value = 33335154
$ws.Range("M100").Value = -33334613

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1107.102
$ws.Range("I112").Value = 1226.6666
$ws.Range("J112").Value = 1099.3043
$ws.Range("K112").Value = 3679.9998
$ws.Range("L112").Value = 3297.9129
$ws.Range("M112").Value = -2571.9998
$ws.Range("N112").Value = -5513.9129

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 842.65515
$ws.Range("J129").Value = 889.8823
$ws.Range("L129").Value = 2669.6469
$ws.Range("N129").Value = -12669.6469

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2677.7632
$ws.Range("I138").Value = 1400.95
$ws.Range("J138").Value = 4096.4443
$ws.Range("K138").Value = 4202.85
$ws.Range("L138").Value = 12289.3329
$ws.Range("M138").Value = 937.1499999999996
$ws.Range("N138").Value = -22569.3329

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 1797.1428
$ws.Range("I141").Value = 1797.1428
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 5391.428400000001
$ws.Range("L141").Value = 0
$ws.Range("M141").Value = -211.4284000000007
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4569.077
$ws.Range("I86").Value = 5657.143
$ws.Range("J86").Value = 3299.6667
$ws.Range("K86").Value = 5657.143
$ws.Range("L86").Value = 3299.6667
$ws.Range("M86").Value = -4534.143
$ws.Range("N86").Value = -5545.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 4569.077
$ws.Range("I89").Value = 5657.143
$ws.Range("J89").Value = 3299.6667
$ws.Range("K89").Value = 28285.715
$ws.Range("L89").Value = 16498.3335
$ws.Range("M89").Value = -22669.715
$ws.Range("N89").Value = -27730.3335

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 631.6875
$ws.Range("I94").Value = 508.91666
$ws.Range("K94").Value = 508.91666
$ws.Range("M94").Value = -57.91665999999998

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1111.1177
$ws.Range("I99").Value = 910.6923
$ws.Range("J99").Value = 1762.5
$ws.Range("K99").Value = 910.6923
$ws.Range("L99").Value = 1762.5
$ws.Range("M99").Value = 587.3077
$ws.Range("N99").Value = -4758.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H140").Value = 42905
$ws.Range("J140").Value = 42905
$ws.Range("L140").Value = 42905
$ws.Range("N140").Value = -53265

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 1139
$ws.Range("I7").Value = 1411.25
$ws.Range("J7").Value = 50
$ws.Range("K7").Value = 1411.25
$ws.Range("L7").Value = 50
$ws.Range("M7").Value = -1298.25
$ws.Range("N7").Value = -276

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1156.8462
$ws.Range("I22").Value = 1467.125
$ws.Range("J22").Value = 660.4
$ws.Range("K22").Value = 1467.125
$ws.Range("L22").Value = 660.4
$ws.Range("M22").Value = -1117.125
$ws.Range("N22").Value = -1360.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2398.2068
$ws.Range("I58").Value = 1998.4286
$ws.Range("J58").Value = 2771.3333
$ws.Range("K58").Value = 1998.4286
$ws.Range("L58").Value = 2771.3333
$ws.Range("M58").Value = -1795.4286
$ws.Range("N58").Value = -3177.3333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2352.1292
$ws.Range("I132").Value = 2102.6191
$ws.Range("J132").Value = 2876.1
$ws.Range("K132").Value = 6307.8573
$ws.Range("L132").Value = 8628.299999999999
$ws.Range("M132").Value = -3777.8573
$ws.Range("N132").Value = -13688.3

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 2398.2068
$ws.Range("I136").Value = 1998.4286
$ws.Range("J136").Value = 2771.3333
$ws.Range("K136").Value = 5995.2858
$ws.Range("L136").Value = 8313.999899999999
$ws.Range("M136").Value = -3445.2858
$ws.Range("N136").Value = -13413.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H140").Value = 75041.336
$ws.Range("J140").Value = 75041.336
$ws.Range("L140").Value = 75041.336
$ws.Range("N140").Value = -85401.336

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 164066.67
$ws.Range("I9").Value = 0
$ws.Range("J9").Value = 164066.67
$ws.Range("K9").Value = 0
$ws.Range("L9").Value = 492200.01
$ws.Range("N9").Value = -492648.01
$ws.Range("M9").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = 0
$ws.Range("K25").Value = 0
$ws.Range("M25").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 8659.333000000001
$ws.Range("J49").Value = 8659.333000000001
$ws.Range("L49").Value = 25977.999
$ws.Range("N49").Value = -26289.999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1654.6666
$ws.Range("I97").Value = 2067.2222
$ws.Range("J97").Value = 1242.1111
$ws.Range("K97").Value = 2067.2222
$ws.Range("L97").Value = 1242.1111
$ws.Range("M97").Value = -1571.2222
$ws.Range("N97").Value = -2234.1111

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 41071.453
$ws.Range("J138").Value = 41071.453
$ws.Range("L138").Value = 41071.453
$ws.Range("N138").Value = -51351.453

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1470.9445
$ws.Range("I82").Value = 1372.4
$ws.Range("J82").Value = 1594.125
$ws.Range("K82").Value = 1372.4
$ws.Range("L82").Value = 1594.125
$ws.Range("M82").Value = -1011.4
$ws.Range("N82").Value = -2316.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1470.9445
$ws.Range("I85").Value = 1372.4
$ws.Range("J85").Value = 1594.125
$ws.Range("K85").Value = 1372.4
$ws.Range("L85").Value = 1594.125
$ws.Range("M85").Value = -124.4000000000001
$ws.Range("N85").Value = -4090.125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 8156888
$ws.Range("I136").Value = 39121.594
$ws.Range("J136").Value = 23812580
$ws.Range("K136").Value = 117364.782
$ws.Range("L136").Value = 71437740
$ws.Range("M136").Value = -114814.782
$ws.Range("N136").Value = -71442840

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2272.6
$ws.Range("I81").Value = 2164.5715
$ws.Range("J81").Value = 2330.7693
$ws.Range("K81").Value = 4329.143
$ws.Range("L81").Value = 4661.5386
$ws.Range("M81").Value = -3268.143
$ws.Range("N81").Value = -6783.5386

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 2272.6
$ws.Range("I84").Value = 2164.5715
$ws.Range("J84").Value = 2330.7693
$ws.Range("K84").Value = 21645.715
$ws.Range("L84").Value = 23307.693
$ws.Range("M84").Value = -16341.715
$ws.Range("N84").Value = -33915.693

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1110
$ws.Range("I96").Value = 890
$ws.Range("J96").Value = 1275
$ws.Range("K96").Value = 890
$ws.Range("L96").Value = 1275
$ws.Range("M96").Value = 483
$ws.Range("N96").Value = -4021

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

Write-Host "Applied market price updates across ALC, BSM, CRP, CUL, GSM, LTW, WVR sheets."
